# Edit filters to enable choosing multiple values and Replace tasks.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new task row (row 14) that was previously blank.
$ws.Range("C14").Value = "Patrik Masrna"
$ws.Range("D14").Value = "UI - Bid progress filters"
$ws.Range("E14").Value = "Edit filters in bid progress to enable choosing multiple values"

# Update the active selection to reflect where the author ended up editing.
$ws.Range("C14").Select()
